$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.089.91"
$ws.Range("E2").Value = "  -1.53%  "

$ws.Range("D3").Value = "'1.992.06"
$ws.Range("E3").Value = "  -2.37%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'240.94"
$ws.Range("E5").Value = "  -6.72%  "

$ws.Range("D6").Value = "'0.602"
$ws.Range("E6").Value = "  -3.63%  "

$ws.Range("D8").Value = "'54.85"
$ws.Range("E8").Value = "  -5.25%  "

$ws.Range("D9").Value = "'0.372"
$ws.Range("E9").Value = "  -4.19%  "

$ws.Range("D10").Value = "'58.70"
$ws.Range("E10").Value = "  +2.63%  "

$ws.Range("D11").Value = "'0.0752"
$ws.Range("E11").Value = "  -6.06%  "

$ws.Range("D12").Value = "'0.0983"
$ws.Range("E12").Value = "  -4.27%  "

$ws.Range("D13").Value = "'2.286.34"
$ws.Range("E13").Value = "  -2.39%  "

$ws.Range("D14").Value = "'14.06"
$ws.Range("E14").Value = "  -5.14%  "

$ws.Range("D15").Value = "'21.12"
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").Value = "'0.756"
$ws.Range("E16").Value = "  -7.97%  "

$ws.Range("E17").Value = "  -6.06%  "

$ws.Range("D18").Value = "'2.004.69"
$ws.Range("E18").Value = "  -2.52%  "

$ws.Range("D19").Value = "'37.009.73"
$ws.Range("E19").Value = "  -1.41%  "

$ws.Range("D20").Value = "'68.11"
$ws.Range("E20").Value = "  -2.87%  "

$ws.Range("E21").Value = "  -5.50%  "

$ws.Range("D22").Value = "'228.53"
$ws.Range("E22").Value = "  -0.56%  "

$ws.Range("E23").Value = "  -4.97%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  -9.56%  "

$ws.Range("D26").Value = "'2.35"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").Value = "'161.33"
$ws.Range("E27").Value = "  -1.46%  "

$ws.Range("E28").Value = "  -5.55%  "

$ws.Range("D29").Value = "'19.08"
$ws.Range("E29").Value = "  -4.93%  "

$ws.Range("D30").Value = "'0.123"
$ws.Range("E30").Value = "  -10.53%  "

$ws.Range("D31").Value = "'1.28"
$ws.Range("E31").Value = "  -5.15%  "

$ws.Range("E32").Value = "  -3.19%  "

$ws.Range("E33").Value = "  -7.11%  "

$ws.Range("E34").Value = "  -7.91%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'4.22"
$ws.Range("E35").Value = "  -6.79%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.34"
$ws.Range("E36").Value = "  -6.28%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "'1.79"
$ws.Range("E38").Value = "  -1.62%  "

$ws.Range("E39").Value = "  -4.42%  "

$ws.Range("D40").Value = "'5.22"
$ws.Range("E40").Value = "  -2.99%  "

$ws.Range("D41").Value = "'3.04"
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").Value = "'1.433.81"
$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0203"
$ws.Range("E43").Value = "  -6.43%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.12"
$ws.Range("E44").Value = "  -5.85%  "

$ws.Range("D45").Value = "'0.0886"
$ws.Range("E45").Value = "  -8.54%  "

$ws.Range("D46").Value = "'88.44"
$ws.Range("E46").Value = "  -3.07%  "

$ws.Range("E47").Value = "  -4.92%  "

$ws.Range("E48").Value = "  -4.87%  "

$ws.Range("D49").Value = "'2.88"
$ws.Range("E49").Value = "  +0.28%  "

$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "'3.64"
$ws.Range("E50").Value = "  +13.05%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'6.62"
$ws.Range("E51").Value = "  -10.75%  "
